$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns C (Summary totals) and D (Total work time) for rows 2-11
$ws.Range("C2").Value = 120097.5
$ws.Range("D2").Value = 334.25

$ws.Range("C3").Value = 58665
$ws.Range("D3").Value = 218

$ws.Range("C4").Value = 52357.5
$ws.Range("D4").Value = 204.5

$ws.Range("C5").Value = 75120
$ws.Range("D5").Value = 206.75

$ws.Range("C6").Value = 97462.5
$ws.Range("D6").Value = 206.5

$ws.Range("C7").Value = 121770
$ws.Range("D7").Value = 206.5

$ws.Range("C8").Value = 146527.5
$ws.Range("D8").Value = 206.5

$ws.Range("C9").Value = 171307.5
$ws.Range("D9").Value = 206.5

$ws.Range("C10").Value = 195930
$ws.Range("D10").Value = 206.5

$ws.Range("C11").Value = 221047.5
$ws.Range("D11").Value = 206.5
